$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 279, shifting the existing rows 279-288 down to 281-290
$ws.Rows.Item(279).Resize(2).Insert()

# Row 279 - new data entry
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value = 44753
$ws.Cells.Item(279, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100101
$ws.Cells.Item(279, 8).Value = "Berries"
$ws.Cells.Item(279, 9).Value = 100101007
$ws.Cells.Item(279, 10).Value = "Kiwi"
$ws.Cells.Item(279, 11).Value = "Hayward"
$ws.Cells.Item(279, 12).Value = "Primera"
$ws.Cells.Item(279, 13).Value = 200
$ws.Cells.Item(279, 14).Value = 14000
$ws.Cells.Item(279, 15).Value = 14000
$ws.Cells.Item(279, 16).Value = 14000
$ws.Cells.Item(279, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(279, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(279, 19).Value = 933
$ws.Cells.Item(279, 20).Value = 15

# Row 280 - new data entry
$ws.Cells.Item(280, 1).Value = 4
$ws.Cells.Item(280, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(280, 3).Value = "Los Lagos"
$ws.Cells.Item(280, 4).Value = 44753
$ws.Cells.Item(280, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(280, 5).Value = 10
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100101
$ws.Cells.Item(280, 8).Value = "Berries"
$ws.Cells.Item(280, 9).Value = 100101007
$ws.Cells.Item(280, 10).Value = "Kiwi"
$ws.Cells.Item(280, 11).Value = "Hayward"
$ws.Cells.Item(280, 12).Value = "Segunda"
$ws.Cells.Item(280, 13).Value = 100
$ws.Cells.Item(280, 14).Value = 12000
$ws.Cells.Item(280, 15).Value = 12000
$ws.Cells.Item(280, 16).Value = 12000
$ws.Cells.Item(280, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(280, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(280, 19).Value = 800
$ws.Cells.Item(280, 20).Value = 15
